$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Bump the start date (B2) by one day; all the other dates in column B are
# formulas relative to B2 (B2+2, B2+7, B4+2, ...) so they recalc automatically.
$ws.Range("B2").Value = 43375

# Mark cm010 (row 11) as added/complete.
$ws.Range("C11").Value = $true

# Update the active selection to B3 (matches the saved cursor position).
$ws.Range("B3").Select()
